$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Epha3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.030023666666667
$ws.Range("H2").Value = 6.090071
$ws.Range("I2").Value = 0.8776223887075381
$ws.Range("J2").Value = 0.8776223887075382
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05620966666666666
$ws.Range("N2").Value = 0.168629
$ws.Range("O2").Value = 0.003368847317172573
$ws.Range("P2").Value = 0.003368847317172572
$ws.Range("Q2").Value = 0.1141069536287778
$ws.Range("R2").Value = 1.026962582659
$ws.Range("S2").Value = 0.002956575829687974
$ws.Range("T2").Value = 0.002956575829687974
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Epha3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.030023666666667
$ws.Range("H3").Value = 6.090071
$ws.Range("I3").Value = 0.8776223887075381
$ws.Range("J3").Value = 0.8776223887075382
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.85444966666667
$ws.Range("N3").Value = 44.563349
$ws.Range("O3").Value = 0.8902805491515402
$ws.Range("P3").Value = 0.8902805491515401
$ws.Range("Q3").Value = 30.15488437864211
$ws.Range("R3").Value = 271.393959407779
$ws.Range("S3").Value = 0.7813301421662334
$ws.Range("T3").Value = 0.7813301421662334
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Epha3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.030023666666667
$ws.Range("H4").Value = 6.090071
$ws.Range("I4").Value = 0.8776223887075381
$ws.Range("J4").Value = 0.8776223887075382
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02237433333333333
$ws.Range("N4").Value = 0.067123
$ws.Range("O4").Value = 0.001340974200585751
$ws.Range("P4").Value = 0.001340974200585751
$ws.Range("Q4").Value = 0.04542042619255555
$ws.Range("R4").Value = 0.408783835733
$ws.Range("S4").Value = 0.001176868981113248
$ws.Range("T4").Value = 0.001176868981113248
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Epha3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.030023666666667
$ws.Range("H5").Value = 6.090071
$ws.Range("I5").Value = 0.8776223887075381
$ws.Range("J5").Value = 0.8776223887075382
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.752099666666667
$ws.Range("N5").Value = 5.256299
$ws.Range("O5").Value = 0.1050096293307016
$ws.Range("P5").Value = 0.1050096293307016
$ws.Range("Q5").Value = 3.556803789692112
$ws.Range("R5").Value = 32.011234107229
$ws.Range("S5").Value = 0.09215880173050349
$ws.Range("T5").Value = 0.09215880173050348
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Epha3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.283071
$ws.Range("H6").Value = 0.849213
$ws.Range("I6").Value = 0.1223776112924619
$ws.Range("J6").Value = 0.1223776112924619
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05620966666666666
$ws.Range("N6").Value = 0.168629
$ws.Range("O6").Value = 0.003368847317172573
$ws.Range("P6").Value = 0.003368847317172572
$ws.Range("Q6").Value = 0.015911326553
$ws.Range("R6").Value = 0.143201938977
$ws.Range("S6").Value = 0.0004122714874845981
$ws.Range("T6").Value = 0.0004122714874845981
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Epha3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.283071
$ws.Range("H7").Value = 0.849213
$ws.Range("I7").Value = 0.1223776112924619
$ws.Range("J7").Value = 0.1223776112924619
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.85444966666667
$ws.Range("N7").Value = 44.563349
$ws.Range("O7").Value = 0.8902805491515402
$ws.Range("P7").Value = 0.8902805491515401
$ws.Range("Q7").Value = 4.204863921593001
$ws.Range("R7").Value = 37.843775294337
$ws.Range("S7").Value = 0.1089504069853067
$ws.Range("T7").Value = 0.1089504069853067
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efna5"
$ws.Range("C8").Value = "Epha3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.283071
$ws.Range("H8").Value = 0.849213
$ws.Range("I8").Value = 0.1223776112924619
$ws.Range("J8").Value = 0.1223776112924619
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02237433333333333
$ws.Range("N8").Value = 0.067123
$ws.Range("O8").Value = 0.001340974200585751
$ws.Range("P8").Value = 0.001340974200585751
$ws.Range("Q8").Value = 0.006333524911
$ws.Range("R8").Value = 0.057001724199
$ws.Range("S8").Value = 0.0001641052194725028
$ws.Range("T8").Value = 0.0001641052194725028
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efna5"
$ws.Range("C9").Value = "Epha3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.283071
$ws.Range("H9").Value = 0.849213
$ws.Range("I9").Value = 0.1223776112924619
$ws.Range("J9").Value = 0.1223776112924619
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 1.752099666666667
$ws.Range("N9").Value = 5.256299
$ws.Range("O9").Value = 0.1050096293307016
$ws.Range("P9").Value = 0.1050096293307016
$ws.Range("Q9").Value = 0.4959686047430001
$ws.Range("R9").Value = 4.463717442687
$ws.Range("S9").Value = 0.0128508276001981
$ws.Range("T9").Value = 0.0128508276001981
